$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 19738.455
$ws.Range("I32").Value = 19622
$ws.Range("J32").Value = 19835.5
$ws.Range("K32").Value = 19622
$ws.Range("L32").Value = 19835.5
$ws.Range("M32").Value = -19296
$ws.Range("N32").Value = -20487.5
$ws.Range("H41").Value = 2037.35
$ws.Range("I41").Value = 1140.8125
$ws.Range("K41").Value = 1140.8125
$ws.Range("M41").Value = -700.8125
$ws.Range("H64").Value = 8861.444
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 8861.444
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = ""
$ws.Range("N64").Value = -9357.444
$ws.Range("H67").Value = 8861.444
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 8861.444
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = ""
$ws.Range("N67").Value = -10577.444
$ws.Range("H112").Value = 1955.25
$ws.Range("I112").Value = 1113.2858
$ws.Range("J112").Value = 3134
$ws.Range("K112").Value = 3339.8574
$ws.Range("L112").Value = 9402
$ws.Range("M112").Value = -2231.8574
$ws.Range("N112").Value = -11618
$ws.Range("H129").Value = 1475.3448
$ws.Range("I129").Value = 1047.1904
$ws.Range("K129").Value = 3141.5712
$ws.Range("M129").Value = 1858.4288
$ws.Range("H137").Value = 3028.6667
$ws.Range("I137").Value = 790.8
$ws.Range("J137").Value = 4147.6
$ws.Range("K137").Value = 2372.4
$ws.Range("L137").Value = 12442.8
$ws.Range("M137").Value = 177.6000000000004
$ws.Range("N137").Value = -17542.8
$ws.Range("H138").Value = 2824.9805
$ws.Range("I138").Value = 1898.0358
$ws.Range("J138").Value = 3953.4348
$ws.Range("K138").Value = 5694.107400000001
$ws.Range("L138").Value = 11860.3044
$ws.Range("M138").Value = -554.1074000000008
$ws.Range("N138").Value = -22140.3044
$ws.Range("H141").Value = 5779.905
$ws.Range("I141").Value = 3119.3
$ws.Range("K141").Value = 9357.900000000001
$ws.Range("M141").Value = -4177.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1297.7778
$ws.Range("I4").Value = 247.16667
$ws.Range("K4").Value = 247.16667
$ws.Range("M4").Value = -131.16667
$ws.Range("H46").Value = 15314.909
$ws.Range("J46").Value = 6484.4287
$ws.Range("L46").Value = 6484.4287
$ws.Range("N46").Value = -7122.4287
$ws.Range("H61").Value = 5457.4
$ws.Range("I61").Value = 5457.4
$ws.Range("K61").Value = 5457.4
$ws.Range("M61").Value = -5245.4
$ws.Range("H97").Value = 1344
$ws.Range("I97").Value = 1257.4117
$ws.Range("J97").Value = 1589.3334
$ws.Range("K97").Value = 1257.4117
$ws.Range("L97").Value = 1589.3334
$ws.Range("M97").Value = -761.4117000000001
$ws.Range("N97").Value = -2581.3334
$ws.Range("H122").Value = 3840.1667
$ws.Range("I122").Value = 2864.889
$ws.Range("K122").Value = 8594.667000000001
$ws.Range("M122").Value = -6144.667000000001
$ws.Range("H132").Value = 4811.5293
$ws.Range("I132").Value = 2271.5715
$ws.Range("K132").Value = 6814.7145
$ws.Range("M132").Value = -4284.7145
$ws.Range("H136").Value = 5457.4
$ws.Range("I136").Value = 5457.4
$ws.Range("K136").Value = 16372.2
$ws.Range("M136").Value = -13822.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("H86").Value = 7145.5454
$ws.Range("I86").Value = 5369
$ws.Range("J86").Value = 8626
$ws.Range("K86").Value = 5369
$ws.Range("L86").Value = 8626
$ws.Range("M86").Value = -4246
$ws.Range("N86").Value = -10872
$ws.Range("H89").Value = 7145.5454
$ws.Range("I89").Value = 5369
$ws.Range("J89").Value = 8626
$ws.Range("K89").Value = 26845
$ws.Range("L89").Value = 43130
$ws.Range("M89").Value = -21229
$ws.Range("N89").Value = -54362
$ws.Range("H134").Value = 2802.9644
$ws.Range("I134").Value = 1379.04
$ws.Range("J134").Value = 14669
$ws.Range("K134").Value = 4137.12
$ws.Range("L134").Value = 44007
$ws.Range("M134").Value = -1602.12
$ws.Range("N134").Value = -49077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8029.1
$ws.Range("I62").Value = 3357.4
$ws.Range("J62").Value = 12700.8
$ws.Range("K62").Value = 3357.4
$ws.Range("L62").Value = 12700.8
$ws.Range("M62").Value = -2733.4
$ws.Range("N62").Value = -13948.8
$ws.Range("H65").Value = 8029.1
$ws.Range("I65").Value = 3357.4
$ws.Range("J65").Value = 12700.8
$ws.Range("K65").Value = 16787
$ws.Range("L65").Value = 63504
$ws.Range("M65").Value = -13667
$ws.Range("N65").Value = -69744

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 17008
$ws.Range("I88").Value = 15000
$ws.Range("J88").Value = 19016
$ws.Range("K88").Value = 45000
$ws.Range("L88").Value = 57048
$ws.Range("H91").Value = 17008
$ws.Range("I91").Value = 15000
$ws.Range("J91").Value = 19016
$ws.Range("K91").Value = 45000
$ws.Range("L91").Value = 57048
$ws.Range("H122").Value = 1341.8276
$ws.Range("I122").Value = 264.33334
$ws.Range("J122").Value = 1466.1538
$ws.Range("K122").Value = 2379.00006
$ws.Range("L122").Value = 13195.3842
$ws.Range("M122").Value = 70.9999399999997
$ws.Range("N122").Value = -18095.3842
$ws.Range("H124").Value = 2000
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H131").Value = 9724316
$ws.Range("J131").Value = 12964746
$ws.Range("L131").Value = 38894238
$ws.Range("N131").Value = -38904318

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2741.2856
$ws.Range("I16").Value = 2932.3333
$ws.Range("J16").Value = 1595
$ws.Range("K16").Value = 2932.3333
$ws.Range("L16").Value = 1595
$ws.Range("M16").Value = -2762.3333
$ws.Range("N16").Value = -1935
$ws.Range("H61").Value = 8059.96
$ws.Range("I61").Value = 6322.5557
$ws.Range("K61").Value = 6322.5557
$ws.Range("M61").Value = -6120.5557
$ws.Range("H68").Value = 3404.6667
$ws.Range("I68").Value = 2542.9062
$ws.Range("J68").Value = 10298.75
$ws.Range("K68").Value = 2542.9062
$ws.Range("L68").Value = 10298.75
$ws.Range("M68").Value = -1793.9062
$ws.Range("N68").Value = -11796.75
$ws.Range("H71").Value = 3404.6667
$ws.Range("I71").Value = 2542.9062
$ws.Range("J71").Value = 10298.75
$ws.Range("K71").Value = 12714.531
$ws.Range("L71").Value = 51493.75
$ws.Range("M71").Value = -8970.530999999999
$ws.Range("N71").Value = -58981.75
$ws.Range("H93").Value = 1501
$ws.Range("I93").Value = 1501
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1501
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = ""
$ws.Range("H113").Value = 8059.96
$ws.Range("I113").Value = 6322.5557
$ws.Range("K113").Value = 6322.5557
$ws.Range("M113").Value = -4152.5557
$ws.Range("H132").Value = 2341.36
$ws.Range("I132").Value = 2341.36
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7024.08
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("H136").Value = 4433.6665
$ws.Range("I136").Value = 2046.4667
$ws.Range("K136").Value = 6139.4001
$ws.Range("M136").Value = -3589.4001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13633
$ws.Range("J45").Value = 13633
$ws.Range("L45").Value = 13633
$ws.Range("N45").Value = -14615
$ws.Range("H132").Value = 2512.6667
$ws.Range("I132").Value = 2431.1777
$ws.Range("J132").Value = 3735
$ws.Range("K132").Value = 7293.533100000001
$ws.Range("L132").Value = 11205
$ws.Range("M132").Value = -4763.533100000001
$ws.Range("N132").Value = -16265
